# fix(MSE-1154): Fix inline styler.
#
# Rename the worksheet, reset the active selection to A1, size the three
# data columns, and make sure the (previously implicit/blank) B2 cell that
# sits under the "Active/Inactive" validation list is materialized with an
# explicit blank value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet was renamed as part of the regeneration of this fixture.
$ws.Name = "24b3b3f3"

# Columns A:C widened to a uniform width (stored column width of 15 in the
# saved XML corresponds to a COM ColumnWidth of 14 + 5/6).
$ws.Range("A1:C1").ColumnWidth = 14.166666666666666

# B2 is part of the "Active,Inactive" list validation (together with C2)
# but previously had no value at all; make sure it is explicitly blank.
$ws.Range("B2").Value = ""

# Selection moves from the stale L6 reference back to A1.
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
